$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "26.742.35"
$ws.Cells.Item(2, 5).Value = "  -0.23%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.639.70"
$ws.Cells.Item(3, 5).Value = "  -0.58%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.20%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "218.04"
$ws.Cells.Item(5, 5).Value = "  +0.64%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  -0.80%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.63%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -0.62%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.08"
$ws.Cells.Item(10, 5).Value = "  -0.53%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.11%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "1.866.69"
$ws.Cells.Item(12, 5).Value = "  -0.66%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.626.61"
$ws.Cells.Item(13, 5).Value = "  -1.16%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  -1.32%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  -1.45%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "64.51"
$ws.Cells.Item(16, 5).Value = "  -1.34%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "26.708.65"
$ws.Cells.Item(17, 5).Value = "  -0.36%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -2.16%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "211.42"
$ws.Cells.Item(19, 5).Value = "  -3.13%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.17%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -0.55%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.19"
$ws.Cells.Item(22, 5).Value = "  -1.11%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "2.31"
$ws.Cells.Item(23, 5).Value = "  -5.25%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "9.26"
$ws.Cells.Item(24, 5).Value = "  -2.49%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "146.85"
$ws.Cells.Item(25, 5).Value = "  +0.21%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.32%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.118"
$ws.Cells.Item(27, 5).Value = "  -1.77%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -0.71%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.57"
$ws.Cells.Item(29, 5).Value = "  -1.04%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -3.13%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.19"
$ws.Cells.Item(31, 5).Value = "  +0.66%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +0.36%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -0.54%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "1.266.01"
$ws.Cells.Item(34, 5).Value = "  -1.32%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -0.87%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -0.78%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -1.93%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.527"
$ws.Cells.Item(38, 5).Value = "  -1.75%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -2.81%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -0.23%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.803"
$ws.Cells.Item(41, 5).Value = "  -1.46%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  -3.02%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "RocketPoolETH"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(43, 4).Value = "1.776.87"
$ws.Cells.Item(43, 5).Value = "  -0.73%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.27"
$ws.Cells.Item(44, 5).Value = "  -3.69%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "91.36"
$ws.Cells.Item(45, 5).Value = "  -0.69%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "60.15"
$ws.Cells.Item(46, 5).Value = "  +0.96%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -1.80%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +0.41%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "7.53"
$ws.Cells.Item(49, 5).Value = "  -3.01%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0960"
$ws.Cells.Item(50, 5).Value = "  -1.01%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "USDD"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.01"
$ws.Cells.Item(51, 5).Value = "  -0.17%  "
